# break out stock.yaml completed
# Fix E73:E75 (bsecode) to be stored as numbers instead of text, matching
# the numeric storage used for bsecode in the rest of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

$ws.Cells.Item(73, 5).Value = 509930
$ws.Cells.Item(74, 5).Value = 590024
$ws.Cells.Item(75, 5).Value = 543220

# Append newly scraped rows 76-78 (chartink screener run at 26/06/2024 11:37:54)
# Column E (bsecode) keeps its original text storage for these fresh rows.

# Row 76: SUPREMEIND
$ws.Cells.Item(76, 1).Value = "26/06/2024 11:37:54"
$ws.Cells.Item(76, 2).Value = 1
$ws.Cells.Item(76, 3).Value = "SUPREMEIND"
$ws.Cells.Item(76, 4).Value = "Supreme Industries Limited"
$ws.Cells.Item(76, 5).NumberFormat = "@"
$ws.Cells.Item(76, 5).Value = "509930"
$ws.Cells.Item(76, 5).NumberFormat = "General"
$ws.Cells.Item(76, 6).Value = -2.04
$ws.Cells.Item(76, 7).Value = 5767.85
$ws.Cells.Item(76, 8).Value = 143838

# Row 77: FACT
$ws.Cells.Item(77, 1).Value = "26/06/2024 11:37:54"
$ws.Cells.Item(77, 2).Value = 2
$ws.Cells.Item(77, 3).Value = "FACT"
$ws.Cells.Item(77, 4).Value = "Fertilizers And Chemicals Travancore Limited"
$ws.Cells.Item(77, 5).NumberFormat = "@"
$ws.Cells.Item(77, 5).Value = "590024"
$ws.Cells.Item(77, 5).NumberFormat = "General"
$ws.Cells.Item(77, 6).Value = 0.9
$ws.Cells.Item(77, 7).Value = 1008.05
$ws.Cells.Item(77, 8).Value = 2497791

# Row 78: MAXHEALTH
$ws.Cells.Item(78, 1).Value = "26/06/2024 11:37:54"
$ws.Cells.Item(78, 2).Value = 3
$ws.Cells.Item(78, 3).Value = "MAXHEALTH"
$ws.Cells.Item(78, 4).Value = "Max Healthcare Institute Ltd"
$ws.Cells.Item(78, 5).NumberFormat = "@"
$ws.Cells.Item(78, 5).Value = "543220"
$ws.Cells.Item(78, 5).NumberFormat = "General"
$ws.Cells.Item(78, 6).Value = -1.57
$ws.Cells.Item(78, 7).Value = 879.45
$ws.Cells.Item(78, 8).Value = 1893149
